# Commit: "add not important files"
# This test fixture workbook was reworked from a "Demand in the same city but
# with different time windows" test case into a "Split order (city) to two
# with zero distance between them" test case: the city "Tambov" is split
# into "Tambov1" and "Tambov2", a new "Huge demand" vehicle capacity is
# bumped, Orders quantities/time-windows are re-tuned, and the Routes sheet
# gains the extra city-pair distances/times (including a zero-distance
# Tambov1<->Tambov2 hop formatted in scientific notation).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# TestDescription sheet: update the scenario description text.
# ---------------------------------------------------------------------
$wsDesc = $wb.Worksheets.Item("TestDescription")
$wsDesc.Range("A1").Value = "Split order (city) to two with zero distnance between them"

# ---------------------------------------------------------------------
# Cities sheet: Tambov -> Tambov1 / Tambov2, Voronezh shifts down a row.
# ---------------------------------------------------------------------
$wsCities = $wb.Worksheets.Item("Cities")
$wsCities.Range("A3").Value = "Tambov1"
$wsCities.Range("A5").Value = "Voronezh"
$wsCities.Range("A4").Value = "Tambov2"

# ---------------------------------------------------------------------
# Vehicles sheet: bump the "Huge demand" vehicle's capacity.
# ---------------------------------------------------------------------
$wsVehicles = $wb.Worksheets.Item("Vehicles")
$wsVehicles.Range("B2").Value = 130

# ---------------------------------------------------------------------
# Orders sheet: retune demand/time-window rows for the split cities.
# ---------------------------------------------------------------------
$wsOrders = $wb.Worksheets.Item("Orders")

$wsOrders.Range("A2").Value = "Tambov1"
$wsOrders.Range("B2").Value = 80
$wsOrders.Range("C2").Value = 750
$wsOrders.Range("D2").Value = 3
$wsOrders.Range("E2").Value = 13

$wsOrders.Range("A3").Value = "Tambov2"
$wsOrders.Range("B3").Value = 70
$wsOrders.Range("C3").Value = 850
$wsOrders.Range("D3").Value = 15
$wsOrders.Range("E3").Value = 17

$wsOrders.Range("A4").Value = "Voronezh"
$wsOrders.Range("B4").Value = 70
$wsOrders.Range("C4").Value = 1100
$wsOrders.Range("D4").Value = 17
$wsOrders.Range("E4").Value = 24

# ---------------------------------------------------------------------
# Routes sheet: rename existing Tambov legs to Tambov1, then append the
# legs to/from Tambov2 (mirroring the Moscow/Voronezh <-> Tambov1 legs)
# plus the zero-distance Tambov1<->Tambov2 leg.
# ---------------------------------------------------------------------
$wsRoutes = $wb.Worksheets.Item("Routes")

$wsRoutes.Range("B2").Value = "Tambov1"
$wsRoutes.Range("A4").Value = "Tambov1"
$wsRoutes.Range("A5").Value = "Tambov1"
$wsRoutes.Range("B6").Value = "Tambov1"

$wsRoutes.Range("A8").Value = "Moscow"
$wsRoutes.Range("B8").Value = "Tambov2"
$wsRoutes.Range("C8").Value = 3
$wsRoutes.Range("D8").Value = 400

$wsRoutes.Range("A9").Value = "Moscow"
$wsRoutes.Range("B9").Value = "Voronezh"
$wsRoutes.Range("C9").Value = 5
$wsRoutes.Range("D9").Value = 450

$wsRoutes.Range("A10").Value = "Tambov2"
$wsRoutes.Range("B10").Value = "Moscow"
$wsRoutes.Range("C10").Value = 3
$wsRoutes.Range("D10").Value = 400

$wsRoutes.Range("A11").Value = "Tambov2"
$wsRoutes.Range("B11").Value = "Voronezh"
$wsRoutes.Range("C11").Value = 2.1111111111111112
$wsRoutes.Range("D11").Value = 190

$wsRoutes.Range("A12").Value = "Voronezh"
$wsRoutes.Range("B12").Value = "Tambov2"
$wsRoutes.Range("C12").Value = 2.1111111111111112
$wsRoutes.Range("D12").Value = 190

$wsRoutes.Range("A13").Value = "Voronezh"
$wsRoutes.Range("B13").Value = "Moscow"
$wsRoutes.Range("C13").Value = 5
$wsRoutes.Range("D13").Value = 450

$wsRoutes.Range("A14").Value = "Tambov1"
$wsRoutes.Range("B14").Value = "Tambov2"
$wsRoutes.Range("C14").Value = 0
$wsRoutes.Range("C14").NumberFormat = "0.00E+00"
$wsRoutes.Range("D14").Value = 0

$wsRoutes.Range("A15").Value = "Tambov2"
$wsRoutes.Range("B15").Value = "Tambov1"
$wsRoutes.Range("C15").Value = 0
$wsRoutes.Range("C15").NumberFormat = "0.00E+00"
$wsRoutes.Range("D15").Value = 0

# ---------------------------------------------------------------------
# Window/selection bookkeeping: Routes becomes the active/selected sheet
# (it was Orders before), with a fresh selection on each touched sheet.
# ---------------------------------------------------------------------
$wsDesc.Select()
$wsOrders.Range("H18").Select()
$wsRoutes.Select()
$wsRoutes.Range("C18").Select()
